$wb = $excel.ActiveWorkbook

# Sheets 1,2,3,5 use a plain "2050" label in E1; sheet 4 uses a "2041-2050" range label.
$plainLabelSheets = @(1, 2, 3, 5)
$rangeLabelSheet = 4
$totalRowSheets = @(1, 2, 3, 4)   # sheets with a "Total" row at row 13 (A1:E13 -> A1:E12)
$costSheet = 6                     # "Custo Total" sheet with a "Total" row at row 4 (A1:B4 -> A1:B3)

foreach ($idx in $plainLabelSheets) {
    $ws = $wb.Worksheets.Item($idx)
    $cell = $ws.Cells.Item(1, 5)
    $cell.Value = "2050"
}

$ws4 = $wb.Worksheets.Item($rangeLabelSheet)
$cell4 = $ws4.Cells.Item(1, 5)
$cell4.Value = "2041-2050"

foreach ($idx in $totalRowSheets) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Rows.Item(13).Delete()
}

$wsCost = $wb.Worksheets.Item($costSheet)
$wsCost.Rows.Item(4).Delete()
